$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "Заготовки" (blanks list) between "Изделия" and
#    "Платежи", and populate it with the list of blank/material names.
# ---------------------------------------------------------------------------
$izdelia = $wb.Worksheets.Item("Изделия")
$zagotovki = $wb.Worksheets.Add($null, $izdelia)
$zagotovki.Name = "Заготовки"

$blanks = @(
    "Наименование",
    "М-30 большая",
    "М-30 пузатая бородатая",
    "М-20 большая",
    "М-20 пузатая",
    "М-15 большая",
    "М-15 пузатая",
    "М-10 большая",
    "М-10 пузатая",
    "М-5 18 см",
    "М-5 14 см",
    "М-5 12 см",
    "М-5 9см",
    "М-3 14 см",
    " М-3 9 см",
    "Колокольчик малый",
    "Колокольчик большой",
    "Яйцо",
    "Шар малый",
    "Шар большой",
    "Ангел",
    "Неваляшка",
    "Футляр 0,05",
    "куколка сарафан",
    "Яблоко",
    "Груша",
    "Браслет 1",
    "Браслет 2",
    "Браслет 3",
    "Браслет 4",
    "Браслет 5",
    "Браслет 6",
    "--Выберите заготовку--"
)

for ($i = 0; $i -lt $blanks.Length; $i++) {
    $zagotovki.Cells.Item($i + 1, 1).Value = $blanks[$i]
}

# ---------------------------------------------------------------------------
# 2. Fix up the "Платежи" sheet rows (download/import fix: correct dates,
#    categories, prices and comments for the last three payment rows).
# ---------------------------------------------------------------------------
$platezhi = $wb.Worksheets.Item("Платежи")

$platezhi.Cells.Item(8, 1).Value = "Алексеева Анастасия"
$platezhi.Cells.Item(8, 2).Value = "2011-Sep-23 / 12:09"
$platezhi.Cells.Item(8, 3).Value = "15 / пузатая / лицо"
$platezhi.Cells.Item(8, 5).Value = 600
$platezhi.Cells.Item(8, 6).Value = "фы"

$platezhi.Cells.Item(9, 2).Value = "2011-Sep-27 / 21:09"
$platezhi.Cells.Item(9, 3).Value = "15 / большая / картинка"
$platezhi.Cells.Item(9, 5).Value = 1750
$platezhi.Cells.Item(9, 6).Value = "первый"

$platezhi.Cells.Item(10, 1).Value = "Алексеева Анастасия"
$platezhi.Cells.Item(10, 2).Value = "2011-Sep-27 / 21:09"
$platezhi.Cells.Item(10, 3).Value = "ангел  / ангел  / ангел "
$platezhi.Cells.Item(10, 5).Value = 100
$platezhi.Cells.Item(10, 6).Value = "второй"

# ---------------------------------------------------------------------------
# 3. Make "Платежи" the active tab (now the 4th sheet after the insert).
# ---------------------------------------------------------------------------
$platezhi.Activate()
